# Scheduled-runner update: refresh market-price-derived profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the Leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 658.5714
$ws.Range("I4").Value = 403.66666
$ws.Range("J4").Value = 849.75
$ws.Range("K4").Value = 403.66666
$ws.Range("L4").Value = 849.75
$ws.Range("M4").Value = -289.66666
$ws.Range("N4").Value = -1077.75
$ws.Range("H62").Value = 4337.75
$ws.Range("I62").Value = 5053.3335
$ws.Range("J62").Value = 3145.111
$ws.Range("K62").Value = 5053.3335
$ws.Range("L62").Value = 3145.111
$ws.Range("M62").Value = -4429.3335
$ws.Range("N62").Value = -4393.111
$ws.Range("H65").Value = 4337.75
$ws.Range("I65").Value = 5053.3335
$ws.Range("J65").Value = 3145.111
$ws.Range("K65").Value = 25266.6675
$ws.Range("L65").Value = 15725.555
$ws.Range("M65").Value = -22146.6675
$ws.Range("N65").Value = -21965.555
$ws.Range("H115").Value = 3842.5
$ws.Range("I115").Value = 3842.5
$ws.Range("K115").Value = 11527.5
$ws.Range("M115").Value = -9960.5
$ws.Range("H129").Value = 1056.4918
$ws.Range("J129").Value = 977.0192
$ws.Range("L129").Value = 2931.0576
$ws.Range("N129").Value = -12931.0576
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2544.3572
$ws.Range("I2").Value = 2468.4167
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 2468.4167
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = -2355.4167
$ws.Range("N2").Value = -3226
$ws.Range("H45").Value = 2312
$ws.Range("I45").Value = 2253.6155
$ws.Range("J45").Value = 2420.4285
$ws.Range("K45").Value = 2253.6155
$ws.Range("L45").Value = 2420.4285
$ws.Range("M45").Value = -1876.6155
$ws.Range("N45").Value = -3174.4285
$ws.Range("H61").Value = 1179.6842
$ws.Range("I61").Value = 892.7931
$ws.Range("J61").Value = 2104.111
$ws.Range("K61").Value = 892.7931
$ws.Range("L61").Value = 2104.111
$ws.Range("M61").Value = -680.7931
$ws.Range("N61").Value = -2528.111
$ws.Range("H116").Value = 2544.3572
$ws.Range("I116").Value = 2468.4167
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 2468.4167
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = -174.4167000000002
$ws.Range("N116").Value = -7588
$ws.Range("H122").Value = 1831.3334
$ws.Range("I122").Value = 1831.3334
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5494.0002
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3044.0002
$ws.Range("N122").Value = $null
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = $null
$ws.Range("H136").Value = 1179.6842
$ws.Range("I136").Value = 892.7931
$ws.Range("J136").Value = 2104.111
$ws.Range("K136").Value = 2678.3793
$ws.Range("L136").Value = 6312.333
$ws.Range("M136").Value = -128.3793000000001
$ws.Range("N136").Value = -11412.333
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2544.3572
$ws.Range("I3").Value = 2468.4167
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 2468.4167
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = -2354.4167
$ws.Range("N3").Value = -3228
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").Value = $null
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2954.78
$ws.Range("I31").Value = 1615.0834
$ws.Range("J31").Value = 3377.842
$ws.Range("K31").Value = 1615.0834
$ws.Range("L31").Value = 3377.842
$ws.Range("M31").Value = -1320.0834
$ws.Range("N31").Value = -3967.842
$ws.Range("H34").Value = 2954.78
$ws.Range("I34").Value = 1615.0834
$ws.Range("J34").Value = 3377.842
$ws.Range("K34").Value = 1615.0834
$ws.Range("L34").Value = 3377.842
$ws.Range("M34").Value = -1413.0834
$ws.Range("N34").Value = -3781.842
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 91999.664
$ws.Range("J106").Value = 7999.5
$ws.Range("L106").Value = 23998.5
$ws.Range("N106").Value = -25890.5
$ws.Range("H131").Value = 971.55
$ws.Range("I131").Value = 526
$ws.Range("J131").Value = 995
$ws.Range("K131").Value = 1578
$ws.Range("L131").Value = 2985
$ws.Range("M131").Value = 3462
$ws.Range("N131").Value = -13065
$ws.Range("H141").Value = 111115940
$ws.Range("I141").Value = 333336400
$ws.Range("J141").Value = 5716.5
$ws.Range("K141").Value = 1000009200
$ws.Range("L141").Value = 17149.5
$ws.Range("M141").Value = -1000004020
$ws.Range("N141").Value = -27509.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1360.421
$ws.Range("I113").Value = 1349.6154
$ws.Range("J113").Value = 1383.8334
$ws.Range("K113").Value = 1349.6154
$ws.Range("L113").Value = 1383.8334
$ws.Range("M113").Value = 820.3846000000001
$ws.Range("N113").Value = -5723.8334
$ws.Range("H122").Value = 1300
$ws.Range("I122").Value = 1380
$ws.Range("J122").Value = 1100
$ws.Range("K122").Value = 4140
$ws.Range("L122").Value = 3300
$ws.Range("M122").Value = -1690
$ws.Range("N122").Value = -8200
$ws.Range("H132").Value = 2545.3171
$ws.Range("I132").Value = 1770.76
$ws.Range("J132").Value = 3755.5625
$ws.Range("K132").Value = 5312.28
$ws.Range("L132").Value = 11266.6875
$ws.Range("M132").Value = -2782.28
$ws.Range("N132").Value = -16326.6875
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2256.6365
$ws.Range("I7").Value = 2313.6667
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 2313.6667
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = -2201.6667
$ws.Range("N7").Value = -2224
$ws.Range("H93").Value = 2178.8572
$ws.Range("I93").Value = 2750
$ws.Range("J93").Value = 2083.6667
$ws.Range("K93").Value = 2750
$ws.Range("L93").Value = 2083.6667
$ws.Range("M93").Value = -1502
$ws.Range("N93").Value = -4579.6667
$ws.Range("H122").Value = 85474.336
$ws.Range("I122").Value = 113232.445
$ws.Range("J122").Value = 2200
$ws.Range("K122").Value = 339697.335
$ws.Range("L122").Value = 6600
$ws.Range("M122").Value = -337247.335
$ws.Range("N122").Value = -11500
$ws.Range("H126").Value = 2256.6365
$ws.Range("I126").Value = 2313.6667
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 6941.000100000001
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -4471.000100000001
$ws.Range("N126").Value = -10940
